$wb = $excel.ActiveWorkbook

$wsReferencing = $wb.Worksheets.Item("Referencing")
$wsTables = $wb.Worksheets.Item("Tables")

# Referencing!A1 used to hold the string "Named reference" (which is what
# Named_reference / INDIRECT("naMed_reFerence") / INDIRECT(B11&"_"&C11)
# resolved to). Replace it with a formula so that those dependents now
# resolve to a number instead, pruning the now-unused shared string.
$wsReferencing.Range("A1").Formula = "=C4"

# Tables sheet gains a new row 1 / cell A1 that pulls the (now numeric)
# value through from the Referencing sheet.
$wsTables.Range("A1").Formula = "=Referencing!D11"

# Update the selections left behind on each sheet.
$wsReferencing.Range("A2").Select()

# Activate the Tables sheet last so it becomes the selected/active tab.
$wsTables.Activate()
$wsTables.Range("A2").Select()
